$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add rows 2-7 in column A with text values "1" through "6"
# (stored as text / shared strings, matching the "Qty" style in column A)
$ws.Range("A2").Value = "1"
$ws.Range("A3").Value = "2"
$ws.Range("A4").Value = "3"
$ws.Range("A5").Value = "4"
$ws.Range("A6").Value = "5"
$ws.Range("A7").Value = "6"

# Move the selection as seen in the saved workbook
$ws.Range("B11").Select()
